$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.707.62"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "1.925.82"
$ws.Range("E3").Value = "  -1.16%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9978"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "334.63"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9968"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4669"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.41%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4158"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.69%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "48.30"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.07%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.08056"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.09%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.024"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.14%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "22.34"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.82%  "

$ws.Range("D13").Value = "1.931.20"
$ws.Range("E13").Value = "  -1.46%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.012"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.23%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.189"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.52%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "89.92"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.07%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.9981"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001036"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.01%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06592"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.30%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.82"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.22%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").Value = "29.674.37"
$ws.Range("E22").Value = "  -0.19%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.548"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.54"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.51%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.205"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.87%  "

$ws.Range("D26").Value = "2.163.51"
$ws.Range("E26").Value = "  -0.99%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "157.21"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.52%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.95"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.29%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.168"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.670"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "117.94"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.14%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.039"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.01%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.09461"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.95%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.440"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.28%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.450"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.91%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.536"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.08%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.06153"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("E38").Value = "  -1.99%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.471"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.34%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.180"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5926"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.52%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9964"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.45%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "10.28"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.13%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1841"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.72%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.372"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.75%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.242"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.50%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.07557"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.91%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.5594"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.15%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "12.21"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.942"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.43%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "112.80"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.20%  "
